$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.369.41'
$ws.Cells.Item(2, 5).Value = '  +0.07%  '
$ws.Cells.Item(3, 4).Value = '2.068.02'
$ws.Cells.Item(3, 5).Value = '  +0.37%  '
$ws.Cells.Item(4, 5).Value = '  -0.08%  '
$ws.Cells.Item(5, 4).Value = '''234.08'
$ws.Cells.Item(5, 5).Value = '  -0.80%  '
$ws.Cells.Item(6, 5).Value = '  +0.57%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 4).Value = '''56.84'
$ws.Cells.Item(8, 5).Value = '  -1.76%  '
$ws.Cells.Item(9, 5).Value = '  +0.49%  '
$ws.Cells.Item(10, 4).Value = '''0.0763'
$ws.Cells.Item(10, 5).Value = '  +0.53%  '
$ws.Cells.Item(11, 5).Value = '  +0.82%  '
$ws.Cells.Item(12, 4).Value = '2.372.75'
$ws.Cells.Item(12, 5).Value = '  +0.26%  '
$ws.Cells.Item(13, 4).Value = '''14.45'
$ws.Cells.Item(13, 5).Value = '  +1.34%  '
$ws.Cells.Item(14, 4).Value = '''20.80'
$ws.Cells.Item(14, 5).Value = '  -0.14%  '
$ws.Cells.Item(15, 4).Value = '''0.777'
$ws.Cells.Item(15, 5).Value = '  +0.23%  '
$ws.Cells.Item(16, 4).Value = '''5.13'
$ws.Cells.Item(16, 5).Value = '  -0.91%  '
$ws.Cells.Item(17, 4).Value = '2.069.16'
$ws.Cells.Item(17, 5).Value = '  +0.35%  '
$ws.Cells.Item(18, 4).Value = '37.318.42'
$ws.Cells.Item(18, 5).Value = '  -0.63%  '
$ws.Cells.Item(19, 4).Value = '''6.41'
$ws.Cells.Item(19, 5).Value = '  +5.49%  '
$ws.Cells.Item(20, 4).Value = '''69.57'
$ws.Cells.Item(20, 5).Value = '  +1.62%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0813'
$ws.Cells.Item(21, 5).Value = '  +0.66%  '
$ws.Cells.Item(22, 4).Value = '''226.17'
$ws.Cells.Item(22, 5).Value = '  +0.70%  '
$ws.Cells.Item(23, 5).Value = '  -0.07%  '
$ws.Cells.Item(24, 4).Value = '''2.43'
$ws.Cells.Item(24, 5).Value = '  +1.27%  '
$ws.Cells.Item(25, 5).Value = '  -1.85%  '
$ws.Cells.Item(26, 4).Value = '''166.85'
$ws.Cells.Item(26, 5).Value = '  +2.52%  '
$ws.Cells.Item(27, 4).Value = '''8.79'
$ws.Cells.Item(27, 5).Value = '  -0.33%  '
$ws.Cells.Item(28, 5).Value = '  +3.37%  '
$ws.Cells.Item(29, 4).Value = '''18.96'
$ws.Cells.Item(29, 5).Value = '  -1.25%  '
$ws.Cells.Item(30, 5).Value = '  -1.88%  '
$ws.Cells.Item(31, 5).Value = '  -0.73%  '
$ws.Cells.Item(32, 4).Value = '''4.48'
$ws.Cells.Item(32, 5).Value = '  +0.72%  '
$ws.Cells.Item(33, 4).Value = '''0.0619'
$ws.Cells.Item(33, 5).Value = '  -1.35%  '
$ws.Cells.Item(34, 5).Value = '  +3.76%  '
$ws.Cells.Item(35, 4).Value = '''2.49'
$ws.Cells.Item(35, 5).Value = '  -5.73%  '
$ws.Cells.Item(36, 5).Value = '  -0.11%  '
$ws.Cells.Item(37, 5).Value = '  -1.82%  '
$ws.Cells.Item(38, 5).Value = '  -3.19%  '
$ws.Cells.Item(39, 5).Value = '  -2.98%  '
$ws.Cells.Item(40, 5).Value = '  -1.34%  '
$ws.Cells.Item(41, 4).Value = '1.476.33'
$ws.Cells.Item(41, 5).Value = '  +0.49%  '
$ws.Cells.Item(42, 4).Value = '''96.09'
$ws.Cells.Item(42, 5).Value = '  +0.99%  '
$ws.Cells.Item(43, 4).Value = '''0.0935'
$ws.Cells.Item(43, 5).Value = '  -1.33%  '
$ws.Cells.Item(44, 4).Value = '''1.18'
$ws.Cells.Item(44, 5).Value = '  +4.29%  '
$ws.Cells.Item(45, 4).Value = '''0.0212'
$ws.Cells.Item(45, 5).Value = '  +0.56%  '
$ws.Cells.Item(46, 4).Value = '''4.23'
$ws.Cells.Item(46, 5).Value = '  -3.78%  '
$ws.Cells.Item(47, 5).Value = '  +0.23%  '
$ws.Cells.Item(48, 4).Value = '''15.25'
$ws.Cells.Item(48, 5).Value = '  -5.47%  '
$ws.Cells.Item(49, 4).Value = '''7.19'
$ws.Cells.Item(49, 5).Value = '  -0.94%  '
$ws.Cells.Item(50, 5).Value = '  +1.63%  '
$ws.Cells.Item(51, 4).Value = '2.261.21'
$ws.Cells.Item(51, 5).Value = '  +0.29%  '
